# Auto-generated Excel COM-interop script applying scheduled-runner price refresh
# to the Zodiark_Profits Leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1797.8823
$ws.Range("I28").Value = 1562.9166
$ws.Range("K28").Value = 1562.9166
$ws.Range("M28").Value = -1077.9166
$ws.Range("H33").Value = 50000276
$ws.Range("I33").Value = 99.42856999999999
$ws.Range("K33").Value = 99.42856999999999
$ws.Range("M33").Value = 129.57143
$ws.Range("H52").Value = 596.3333
$ws.Range("I52").Value = 596.3333
$ws.Range("K52").Value = 1788.9999
$ws.Range("M52").Value = -1628.9999
$ws.Range("H111").Value = 2215.6
$ws.Range("I111").Value = 1592.6666
$ws.Range("J111").Value = 3150
$ws.Range("K111").Value = 4777.9998
$ws.Range("L111").Value = 9450
$ws.Range("M111").Value = -1710.9998
$ws.Range("N111").Value = -15584
$ws.Range("H116").Value = 2944.0625
$ws.Range("I116").Value = 2944.6428
$ws.Range("K116").Value = 2944.6428
$ws.Range("M116").Value = 497.3571999999999
$ws.Range("H132").Value = 1456.909
$ws.Range("I132").Value = 1325.9667
$ws.Range("K132").Value = 3977.9001
$ws.Range("M132").Value = -1447.9001
$ws.Range("H137").Value = 2602.6843
$ws.Range("I137").Value = 3553.3333
$ws.Range("K137").Value = 10659.9999
$ws.Range("M137").Value = -8109.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = 0
$ws.Range("H5").Value = 101
$ws.Range("I5").Value = 101
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 101
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = 11
$ws.Range("H6").Value = 26250
$ws.Range("I6").Value = 26250
$ws.Range("K6").Value = 26250
$ws.Range("M6").Value = -26077
$ws.Range("H45").Value = 2566.238
$ws.Range("I45").Value = 2499.6924
$ws.Range("J45").Value = 2674.375
$ws.Range("K45").Value = 2499.6924
$ws.Range("L45").Value = 2674.375
$ws.Range("M45").Value = -2122.6924
$ws.Range("N45").Value = -3428.375
$ws.Range("H102").Value = 254999
$ws.Range("I102").Value = 9999
$ws.Range("K102").Value = 9999
$ws.Range("M102").Value = -8377
$ws.Range("H122").Value = 6687.8335
$ws.Range("J122").Value = 5104.6
$ws.Range("L122").Value = 15313.8
$ws.Range("N122").Value = -20213.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 101
$ws.Range("I4").Value = 101
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 101
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = 14
$ws.Range("H80").Value = 735
$ws.Range("I80").Value = 981
$ws.Range("J80").Value = 623.1818
$ws.Range("K80").Value = 981
$ws.Range("L80").Value = 623.1818
$ws.Range("M80").Value = 17
$ws.Range("N80").Value = -2619.1818
$ws.Range("H83").Value = 735
$ws.Range("I83").Value = 981
$ws.Range("J83").Value = 623.1818
$ws.Range("K83").Value = 4905
$ws.Range("L83").Value = 3115.909
$ws.Range("M83").Value = 87
$ws.Range("N83").Value = -13099.909
$ws.Range("H96").Value = 2618.3333
$ws.Range("I96").Value = 2618.3333
$ws.Range("K96").Value = 2618.3333
$ws.Range("M96").Value = 127.6667000000002
$ws.Range("H105").Value = 3079.1667
$ws.Range("I105").Value = 2494.3333
$ws.Range("K105").Value = 2494.3333
$ws.Range("M105").Value = -747.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1525
$ws.Range("J4").Value = 1787.5
$ws.Range("L4").Value = 1787.5
$ws.Range("N4").Value = -2011.5
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 5000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = -4887
$ws.Range("N6").Value = 0
$ws.Range("H7").Value = 447.85715
$ws.Range("I7").Value = 603.5
$ws.Range("K7").Value = 603.5
$ws.Range("M7").Value = -490.5
$ws.Range("H16").Value = 3290.7856
$ws.Range("I16").Value = 3737.75
$ws.Range("J16").Value = 2694.8333
$ws.Range("K16").Value = 3737.75
$ws.Range("L16").Value = 2694.8333
$ws.Range("M16").Value = -3450.75
$ws.Range("N16").Value = -3268.8333
$ws.Range("H23").Value = 10599.8
$ws.Range("I23").Value = 6504.5
$ws.Range("J23").Value = 13330
$ws.Range("K23").Value = 6504.5
$ws.Range("L23").Value = 13330
$ws.Range("M23").Value = -6264.5
$ws.Range("N23").Value = -13810
$ws.Range("H27").Value = 10599.8
$ws.Range("I27").Value = 6504.5
$ws.Range("J27").Value = 13330
$ws.Range("K27").Value = 6504.5
$ws.Range("L27").Value = 13330
$ws.Range("M27").Value = -6312.5
$ws.Range("N27").Value = -13714
$ws.Range("H105").Value = 20973.334
$ws.Range("I105").Value = 36292
$ws.Range("J105").Value = 1825
$ws.Range("K105").Value = 36292
$ws.Range("L105").Value = 1825
$ws.Range("M105").Value = -34545
$ws.Range("N105").Value = -5319
$ws.Range("H113").Value = 3290.7856
$ws.Range("I113").Value = 3737.75
$ws.Range("J113").Value = 2694.8333
$ws.Range("K113").Value = 3737.75
$ws.Range("L113").Value = 2694.8333
$ws.Range("M113").Value = -1567.75
$ws.Range("N113").Value = -7034.8333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 142857520
$ws.Range("H8").Value = 652.25
$ws.Range("I8").Value = 652.25
$ws.Range("K8").Value = 1956.75
$ws.Range("M8").Value = -1817.75
$ws.Range("H38").Value = 175.67857
$ws.Range("J38").Value = 193.9
$ws.Range("L38").Value = 581.7
$ws.Range("N38").Value = -1275.7
$ws.Range("H107").Value = 749
$ws.Range("J107").Value = 911.75
$ws.Range("L107").Value = 2735.25
$ws.Range("N107").Value = -6575.25
$ws.Range("H129").Value = 3408.6667
$ws.Range("I129").Value = 4119
$ws.Range("K129").Value = 12357
$ws.Range("M129").Value = -7357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3763.9092
$ws.Range("J22").Value = 3199.8572
$ws.Range("L22").Value = 3199.8572
$ws.Range("N22").Value = -4257.8572
$ws.Range("H31").Value = 4020
$ws.Range("I31").Value = 4020
$ws.Range("K31").Value = 4020
$ws.Range("M31").Value = -3728
$ws.Range("H37").Value = 4020
$ws.Range("I37").Value = 4020
$ws.Range("K37").Value = 4020
$ws.Range("M37").Value = -3743
$ws.Range("H136").Value = 50544
$ws.Range("J136").Value = 50544
$ws.Range("L136").Value = 151632
$ws.Range("N136").Value = -156732

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3785.3333
$ws.Range("J7").Value = 4853.7144
$ws.Range("L7").Value = 4853.7144
$ws.Range("N7").Value = -5077.7144
$ws.Range("H10").Value = 2374.75
$ws.Range("I10").Value = 999.5
$ws.Range("K10").Value = 999.5
$ws.Range("M10").Value = -859.5
$ws.Range("H40").Value = 7747.4165
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864
$ws.Range("H61").Value = 786.3077
$ws.Range("I61").Value = 786.3077
$ws.Range("K61").Value = 786.3077
$ws.Range("M61").Value = -584.3077
$ws.Range("H100").Value = 5456.4287
$ws.Range("J100").Value = 7999.6665
$ws.Range("L100").Value = 7999.6665
$ws.Range("N100").Value = -9081.666499999999
$ws.Range("H101").Value = 101376
$ws.Range("J101").Value = 101376
$ws.Range("L101").Value = 101376
$ws.Range("N101").Value = -107866
$ws.Range("H113").Value = 786.3077
$ws.Range("I113").Value = 786.3077
$ws.Range("K113").Value = 786.3077
$ws.Range("M113").Value = 1383.6923
$ws.Range("H126").Value = 3785.3333
$ws.Range("J126").Value = 4853.7144
$ws.Range("L126").Value = 14561.1432
$ws.Range("N126").Value = -19501.1432
$ws.Range("H132").Value = 3712.4075
$ws.Range("I132").Value = 3345.48
$ws.Range("K132").Value = 10036.44
$ws.Range("M132").Value = -7506.440000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 424.3125
$ws.Range("I107").Value = 401.7
$ws.Range("J107").Value = 462
$ws.Range("K107").Value = 1205.1
$ws.Range("L107").Value = 1386
$ws.Range("M107").Value = 714.9000000000001
$ws.Range("N107").Value = -5226
$ws.Range("H113").Value = 1253.5172
$ws.Range("I113").Value = 788.7222
$ws.Range("J113").Value = 2014.091
$ws.Range("K113").Value = 2366.1666
$ws.Range("L113").Value = 6042.272999999999
$ws.Range("M113").Value = -196.1666
$ws.Range("N113").Value = -10382.273
$ws.Range("H122").Value = 4766.914
$ws.Range("I122").Value = 5321
$ws.Range("J122").Value = 3829.2307
$ws.Range("K122").Value = 15963
$ws.Range("L122").Value = 11487.6921
$ws.Range("M122").Value = -13513
$ws.Range("N122").Value = -16387.6921
$ws.Range("H130").Value = 47500
$ws.Range("J130").Value = 47500
$ws.Range("L130").Value = 47500
$ws.Range("N130").Value = -57540
